# Auto-generated edit script: updates live market-price derived cells
# across 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching
# the scheduled market-data refresh described by the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1155.3846
$ws.Range("I15").Value = 1155.3846
$ws.Range("K15").Value = 3466.1538
$ws.Range("M15").Value = -3297.1538

$ws.Range("H64").Value = 114097.664
$ws.Range("I64").Value = 252350
$ws.Range("J64").Value = 3495.8
$ws.Range("K64").Value = 252350
$ws.Range("L64").Value = 3495.8
$ws.Range("M64").Value = -252102
$ws.Range("N64").Value = -3991.8

$ws.Range("H67").Value = 114097.664
$ws.Range("I67").Value = 252350
$ws.Range("J67").Value = 3495.8
$ws.Range("K67").Value = 252350
$ws.Range("L67").Value = 3495.8
$ws.Range("M67").Value = -251492
$ws.Range("N67").Value = -5211.8

$ws.Range("H70").Value = 1500.125
$ws.Range("J70").Value = 1000.5
$ws.Range("L70").Value = 3001.5
$ws.Range("N70").Value = -3541.5

$ws.Range("H73").Value = 1500.125
$ws.Range("J73").Value = 1000.5
$ws.Range("L73").Value = 3001.5
$ws.Range("N73").Value = -4873.5

$ws.Range("H110").Value = 29000
$ws.Range("J110").Value = 29000
$ws.Range("L110").Value = 29000
$ws.Range("N110").Value = -37180

$ws.Range("H123").Value = 22335.4
$ws.Range("J123").Value = 22335.4
$ws.Range("L123").Value = 22335.4
$ws.Range("N123").Value = -32135.4

$ws.Range("H129").Value = 846.3158
$ws.Range("I129").Value = 553.44446
$ws.Range("J129").Value = 885.65674
$ws.Range("K129").Value = 1660.33338
$ws.Range("L129").Value = 2656.97022
$ws.Range("M129").Value = 3339.66662
$ws.Range("N129").Value = -12656.97022

$ws.Range("H135").Value = 742.4545000000001
$ws.Range("I135").Value = 643.82355
$ws.Range("J135").Value = 1077.8
$ws.Range("K135").Value = 5794.41195
$ws.Range("L135").Value = 9700.199999999999
$ws.Range("M135").Value = -3259.41195
$ws.Range("N135").Value = -14770.2

$ws.Range("H137").Value = 968
$ws.Range("I137").Value = 965.4761999999999
$ws.Range("J137").Value = 994.5
$ws.Range("K137").Value = 2896.4286
$ws.Range("L137").Value = 2983.5
$ws.Range("M137").Value = -346.4285999999997
$ws.Range("N137").Value = -8083.5

$ws.Range("H138").Value = 4142.1553
$ws.Range("I138").Value = 2138.6
$ws.Range("J138").Value = 4841.07
$ws.Range("K138").Value = 6415.799999999999
$ws.Range("L138").Value = 14523.21
$ws.Range("M138").Value = -1275.799999999999
$ws.Range("N138").Value = -24803.21

$ws.Range("H141").Value = 2674.1538
$ws.Range("I141").Value = 2755.3333
$ws.Range("J141").Value = 1700
$ws.Range("K141").Value = 8265.999899999999
$ws.Range("L141").Value = 5100
$ws.Range("M141").Value = -3085.999899999999
$ws.Range("N141").Value = -15460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42842.168
$ws.Range("I2").Value = 1166.9333
$ws.Range("J2").Value = 112300.89
$ws.Range("K2").Value = 1166.9333
$ws.Range("L2").Value = 112300.89
$ws.Range("M2").Value = -1053.9333
$ws.Range("N2").Value = -112526.89

$ws.Range("H113").Value = 36450
$ws.Range("J113").Value = 36450
$ws.Range("L113").Value = 36450
$ws.Range("N113").Value = -45128

$ws.Range("H116").Value = 42842.168
$ws.Range("I116").Value = 1166.9333
$ws.Range("J116").Value = 112300.89
$ws.Range("K116").Value = 1166.9333
$ws.Range("L116").Value = 112300.89
$ws.Range("M116").Value = 1127.0667
$ws.Range("N116").Value = -116888.89

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42842.168
$ws.Range("I3").Value = 1166.9333
$ws.Range("J3").Value = 112300.89
$ws.Range("K3").Value = 1166.9333
$ws.Range("L3").Value = 112300.89
$ws.Range("M3").Value = -1052.9333
$ws.Range("N3").Value = -112528.89

$ws.Range("H107").Value = 30341752
$ws.Range("I107").Value = 41718280
$ws.Range("J107").Value = 4337
$ws.Range("K107").Value = 41718280
$ws.Range("L107").Value = 4337
$ws.Range("M107").Value = -41716360
$ws.Range("N107").Value = -8177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20673.334
$ws.Range("J59").Value = 20673.334
$ws.Range("L59").Value = 20673.334
$ws.Range("N59").Value = -22963.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6099.175
$ws.Range("I5").Value = 999.2069
$ws.Range("J5").Value = 19544.545
$ws.Range("K5").Value = 2997.6207
$ws.Range("L5").Value = 58633.63499999999
$ws.Range("M5").Value = -2885.6207
$ws.Range("N5").Value = -58857.63499999999

$ws.Range("H50").Value = 1867.375
$ws.Range("I50").Value = 2227.8
$ws.Range("J50").Value = 1266.6666
$ws.Range("K50").Value = 6683.400000000001
$ws.Range("L50").Value = 3799.9998
$ws.Range("M50").Value = -6202.400000000001
$ws.Range("N50").Value = -4761.9998

$ws.Range("H53").Value = 1867.375
$ws.Range("I53").Value = 2227.8
$ws.Range("J53").Value = 1266.6666
$ws.Range("K53").Value = 6683.400000000001
$ws.Range("L53").Value = 3799.9998
$ws.Range("M53").Value = -6202.400000000001
$ws.Range("N53").Value = -4761.9998

$ws.Range("H132").Value = 4875.4165
$ws.Range("J132").Value = 6301
$ws.Range("L132").Value = 56709
$ws.Range("N132").Value = -61769

$ws.Range("H135").Value = 6099.175
$ws.Range("I135").Value = 999.2069
$ws.Range("J135").Value = 19544.545
$ws.Range("K135").Value = 8992.8621
$ws.Range("L135").Value = 175900.905
$ws.Range("M135").Value = -6457.8621
$ws.Range("N135").Value = -180970.905

$ws.Range("H136").Value = 2565.3333
$ws.Range("I136").Value = 1998
$ws.Range("K136").Value = 5994
$ws.Range("M136").Value = -894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 4677
$ws.Range("J47").Value = 4677
$ws.Range("L47").Value = 4677
$ws.Range("N47").Value = -5813

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3963.0908
$ws.Range("J68").Value = 4886.75
$ws.Range("L68").Value = 4886.75
$ws.Range("N68").Value = -6384.75

$ws.Range("H69").Value = 29666.666
$ws.Range("J69").Value = 24500
$ws.Range("L69").Value = 24500
$ws.Range("N69").Value = -26122

$ws.Range("H71").Value = 3963.0908
$ws.Range("J71").Value = 4886.75
$ws.Range("L71").Value = 24433.75
$ws.Range("N71").Value = -31921.75

$ws.Range("H72").Value = 29666.666
$ws.Range("J72").Value = 24500
$ws.Range("L72").Value = 73500
$ws.Range("N72").Value = -81612
